$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 88; this pushes the existing rows 88..135
# down to 89..136 and carries formatting from the surrounding rows.
$ws.Rows.Item(88).Insert()

# Populate the newly inserted row 88 with the new record (same underlying
# record as the old row 88, but with an updated date and volume).
$ws.Range("A88").Value  = 1
$ws.Range("B88").Value  = "Agrícola del Norte S.A. de Arica"
$ws.Range("C88").Value  = "Arica y Parinacota"
$ws.Range("D88").Value  = 44827
$ws.Range("E88").Value  = 15
$ws.Range("F88").Value  = 100112036
$ws.Range("G88").Value  = "Caigua"
$ws.Range("H88").Value  = "Sin especificar"
$ws.Range("I88").Value  = "Primera"
$ws.Range("J88").Value  = 140
$ws.Range("K88").Value  = 6000
$ws.Range("L88").Value  = 7000
$ws.Range("M88").Value  = 6500
$ws.Range("N88").Value  = "$/caja 20 kilos"
$ws.Range("O88").Value  = "Región de Arica y Parinacota"
$ws.Range("P88").Value  = 325
$ws.Range("Q88").Value  = 20
$ws.Range("R88").Value  = "Hortaliza"
